$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.169.02'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.589.09'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.55'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.72%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.84%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.245'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.27%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0604'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.99'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.32%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.812.01'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.571.74'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.83%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.47%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.52%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.59'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.159.21'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0723'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.37'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '214.27'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.49%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.23'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.65%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.96'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.43'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.95'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.38%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.96%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.21%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.53%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.417.76'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +8.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.95'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.39%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.67%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.86%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.584'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -4.38%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.823'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.18%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.88%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.948'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -13.36%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.13'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.32%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.763'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.723.85'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '61.06'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '86.84'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.48'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0501'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.96%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0958'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.26%  '
